$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last game row (row 4) entirely
$ws.Rows.Item(4).Delete()

# Row 2 updates
$ws.Range("D2").Value = 1100

# Row 2: drop-in/goalie list cell no longer used -> clear it
$ws.Range("I2").ClearContents()

# Row 3 updates
$ws.Range("C3").Value = 20230907
$ws.Range("D3").Value = 1100
$ws.Range("E3").Value = "Prospects Athletics"
$ws.Range("F3").Value = 5
$ws.Range("G3").Value = "1;"
$ws.Range("I3").Value = "0;"
